$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Insert the new "Chapter 4 Requirements Engineering" block of
#    paragraphs right before the existing "Ethical principles..." one.
# ---------------------------------------------------------------
$anchor = $d.Paragraphs.Item(1)
for ($n = 0; $n -lt 17; $n++) {
    $anchor.Range.InsertParagraphBefore() | Out-Null
}

# --- new paragraph 1 ---
$p = $d.Paragraphs.Item(1)
$p.Range.Text = "Chapter 4 Requirements Engineering"

# --- new paragraph 2 ---
$p = $d.Paragraphs.Item(2)
$p.Range.Text = "Use cases: identify the actors in an interaction and which describe the interaction itself; A set of use cases should describe all possible interactions with the system; UML sequence diagrams may be used to add detail to use-cases by showing the sequence of event processing in the system."
$lbl = $d.Range($p.Range.Start + 0, $p.Range.Start + 10)
$lbl.Font.Underline = 1

# --- new paragraph 3 ---
$p = $d.Paragraphs.Item(3)
$p.Range.Text = "RE processes: Elicitation, Analysis, Validation, Management. Its an iterative activity in which these processes are interleaved. Requirement elicitation: Start, User requirements elicitation, System req. elicitation; Requirements specification: Business requirements specification, User requirements specification, System requirements specification and modeling; Requirements validation: Feasibility study, Prototyping, Reviews."
$lbl = $d.Range($p.Range.Start + 0, $p.Range.Start + 13)
$lbl.Font.Underline = 1
$lbl = $d.Range($p.Range.Start + 129, $p.Range.Start + 153)
$lbl.Font.Underline = 1
$lbl = $d.Range($p.Range.Start + 217, $p.Range.Start + 244)
$lbl.Font.Underline = 1

# --- new paragraph 4 ---
$p = $d.Paragraphs.Item(4)
$p.Range.Text = "Ways of writing a system requirements specification: Natural language, Structured natural language, Design description languages, Graphical notations, Mathematical specifications."
$lbl = $d.Range($p.Range.Start + 0, $p.Range.Start + 52)
$lbl.Font.Underline = 1

# --- new paragraph 5 ---
$p = $d.Paragraphs.Item(5)
$p.Range.Text = "Metrics for nonfunctional requirements: Speed, Size, Ease of use, Reliability, Robustness, Portability"
$lbl = $d.Range($p.Range.Start + 0, $p.Range.Start + 39)
$lbl.Font.Underline = 1

# --- new paragraph 6 ---
$p = $d.Paragraphs.Item(6)
$p.Range.Text = "Non-functional requirements (Product requirements, Organizational requirements, External requirements) define system properties and constraints e.g. reliability, response time and storage requirements, and may be more critical than functional requirements."

# --- new paragraph 7 ---
$p = $d.Paragraphs.Item(7)
$p.Range.Text = "Requirements completeness and consistency: They should include descriptions of all facilities required. There should be no conflicts or contradictions in the descriptions of the system facilities."
$lbl = $d.Range($p.Range.Start + 0, $p.Range.Start + 42)
$lbl.Font.Underline = 1

# --- new paragraph 8 ---
$p = $d.Paragraphs.Item(8)
$p.Range.Text = "System stakeholder types: End users, System managers, System owners, External Stakeholders"
$lbl = $d.Range($p.Range.Start + 0, $p.Range.Start + 25)
$lbl.Font.Underline = 1

# --- new paragraph 9 ---
$p = $d.Paragraphs.Item(9)
$p.Range.Text = "Agile methods: use incremental requirements engineering and may express requirements as user stories which is practical for business systems but not for critical systems."
$lbl = $d.Range($p.Range.Start + 0, $p.Range.Start + 14)
$lbl.Font.Underline = 1

# --- new paragraph 10 ---
$p = $d.Paragraphs.Item(10)
$p.Range.Text = "Functional requirements: Statements of services the system should provide, how the system should react to particular inputs and how the system should behave in particular situations."
$lbl = $d.Range($p.Range.Start + 0, $p.Range.Start + 24)
$lbl.Font.Underline = 1

# --- new paragraph 11 ---
$p = $d.Paragraphs.Item(11)
$p.Range.Text = "Requirements Abstraction: A contract defined its needs in a sufficiently abstract way that a solution is not pre-defined. The requirements must be written so that several contractors can bid for the contract, offering different ways of meeting the client organization’s needs. Once a contract has been awarded, the contractor must write a system definition for the client in more detail so that the client understands and can validate what the software will do. Both of these documents may be called the requirements document for the system."
$lbl = $d.Range($p.Range.Start + 0, $p.Range.Start + 25)
$lbl.Font.Underline = 1

# --- new paragraph 12 ---
$p = $d.Paragraphs.Item(12)
$p.Range.Text = "System requirements: A structured document setting out detailed descriptions of the system’s functions, services and operational constraints. Defines what should be implemented so may be part of a contract between client and contractor."
$lbl = $d.Range($p.Range.Start + 0, $p.Range.Start + 20)
$lbl.Font.Underline = 1

# --- new paragraph 13 ---
$p = $d.Paragraphs.Item(13)

# --- new paragraph 14 ---
$p = $d.Paragraphs.Item(14)
$p.Range.Text = "The process of establishing the services that a customer requires from a system and the constraints under which it operates and is developed."

# --- new paragraph 15 ---
$p = $d.Paragraphs.Item(15)
$p.Range.Text = "The system requirements are the descriptions of the system services and constraints that are generated during the requirements engineering process may range from a high-level abstract statement of a service or of a system constraint to a detailed mathematical functional specification."

# --- new paragraph 16 ---
$p = $d.Paragraphs.Item(16)

# --- new paragraph 17 ---
$p = $d.Paragraphs.Item(17)

# ---------------------------------------------------------------
# 2) Split the run of the (pre-existing) "Ethical principles:" and
#    "Web software engineering:" paragraphs so that the leading
#    label is its own run; underline the "Ethical principles:" label.
# ---------------------------------------------------------------
$epPara = $null
$wsePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt.StartsWith("Ethical principles:")) { $epPara = $d.Paragraphs.Item($i) }
    if ($txt.StartsWith("Web software engineering:")) { $wsePara = $d.Paragraphs.Item($i) }
}

$epLabel = "Ethical principles:"
$epLbl = $d.Range($epPara.Range.Start, $epPara.Range.Start + $epLabel.Length)
$epLbl.Font.Underline = 1

$wseLabel = "Web software engineering:"
$wseLbl = $d.Range($wsePara.Range.Start, $wsePara.Range.Start + $wseLabel.Length)
$wseLbl.Bold = 1
$wseLbl.Bold = 0

# ---------------------------------------------------------------
# 3) Drop the trailing "Requirements engineering: " paragraph,
#    folding its bookmark into the end of the previous (UML) one.
# ---------------------------------------------------------------
$n = $d.Paragraphs.Count
$reLabel = "Requirements engineering: "
$reText = $d.Range($d.Paragraphs.Item($n).Range.Start, $d.Paragraphs.Item($n).Range.Start + $reLabel.Length)
$reText.Delete()
$n2 = $d.Paragraphs.Count
$prev2 = $d.Paragraphs.Item($n2 - 1)
$mark = $d.Range($prev2.Range.End - 1, $prev2.Range.End)
$mark.Delete()

Write-Output "done"
